# Update "想去人数" (want-to-go count) figures that changed between the two
# data pulls. The same underlying events appear on both the "展览" sheet
# (sheet index 1) and the "全部类型" sheet (sheet index 4, offset by one row
# because it has an extra leading row), so both copies need to be updated.

$wb = $excel.ActiveWorkbook

# Sheet 1: "展览"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 12808
$ws1.Range("F10").Value = 12722
$ws1.Range("F11").Value = 270
$ws1.Range("F12").Value = 20
$ws1.Range("F13").Value = 8643
$ws1.Range("F14").Value = 7632
$ws1.Range("F15").Value = 183

# Sheet 4: "全部类型"
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 12808
$ws4.Range("F11").Value = 12722
$ws4.Range("F12").Value = 270
$ws4.Range("F13").Value = 20
$ws4.Range("F14").Value = 8643
$ws4.Range("F15").Value = 7632
$ws4.Range("F16").Value = 183
